$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 19002.309
$ws.Range("J87").Value = 19002.309
$ws.Range("L87").Value = 19002.309
$ws.Range("N87").Value = -21498.309

$ws.Range("H90").Value = 19002.309
$ws.Range("J90").Value = 19002.309
$ws.Range("L90").Value = 57006.927
$ws.Range("N90").Value = -69486.927

$ws.Range("H94").Value = 11266.667
$ws.Range("I94").Value = 10900
$ws.Range("K94").Value = 10900
$ws.Range("M94").Value = -10449

$ws.Range("H97").Value = 1080
$ws.Range("I97").Value = 250
$ws.Range("J97").Value = 1702.5
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 5107.5
$ws.Range("M97").Value = -254
$ws.Range("N97").Value = -6099.5

$ws.Range("H137").Value = 836.65216
$ws.Range("I137").Value = 784.7
$ws.Range("K137").Value = 2354.1
$ws.Range("M137").Value = 195.8999999999996

$ws.Range("H138").Value = 2323.7705
$ws.Range("I138").Value = 983.3200000000001
$ws.Range("J138").Value = 3254.639
$ws.Range("K138").Value = 2949.96
$ws.Range("L138").Value = 9763.917000000001
$ws.Range("M138").Value = 2190.04
$ws.Range("N138").Value = -20043.917

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5906.391
$ws.Range("I61").Value = 5906.391
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5906.391
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -5694.391
$ws.Range("N61").ClearContents()

$ws.Range("H88").Value = 2722.2222
$ws.Range("I88").Value = 2900
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 2900
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -2494
$ws.Range("N88").Value = -3312

$ws.Range("H91").Value = 2722.2222
$ws.Range("I91").Value = 2900
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 2900
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -1496
$ws.Range("N91").Value = -5308

$ws.Range("H122").Value = 1463.9565
$ws.Range("I122").Value = 1526.3846
$ws.Range("J122").Value = 1382.8
$ws.Range("K122").Value = 4579.1538
$ws.Range("L122").Value = 4148.4
$ws.Range("M122").Value = -2129.1538
$ws.Range("N122").Value = -9048.4

$ws.Range("H132").Value = 37254.53
$ws.Range("I132").Value = 46547.695
$ws.Range("J132").Value = 13505.333
$ws.Range("K132").Value = 139643.085
$ws.Range("L132").Value = 40515.999
$ws.Range("M132").Value = -137113.085
$ws.Range("N132").Value = -45575.999

$ws.Range("H136").Value = 5906.391
$ws.Range("I136").Value = 5906.391
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 17719.173
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15169.173
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1774
$ws.Range("I107").Value = 1669.6666
$ws.Range("K107").Value = 1669.6666
$ws.Range("M107").Value = 250.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 46441.1
$ws.Range("I31").Value = 55180.5
$ws.Range("J31").Value = 4492
$ws.Range("K31").Value = 55180.5
$ws.Range("L31").Value = 4492
$ws.Range("M31").Value = -54885.5
$ws.Range("N31").Value = -5082

$ws.Range("H34").Value = 46441.1
$ws.Range("I34").Value = 55180.5
$ws.Range("J34").Value = 4492
$ws.Range("K34").Value = 55180.5
$ws.Range("L34").Value = 4492
$ws.Range("M34").Value = -54978.5
$ws.Range("N34").Value = -4896

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H56").Value = 40000.668
$ws.Range("J56").Value = 40000.668
$ws.Range("L56").Value = 40000.668
$ws.Range("N56").Value = -41690.668

$ws.Range("H58").Value = 1337
$ws.Range("I58").Value = 1204.4
$ws.Range("K58").Value = 1204.4
$ws.Range("M58").Value = -1001.4

$ws.Range("H132").Value = 5247.8184
$ws.Range("I132").Value = 5316.643
$ws.Range("J132").Value = 4862.4
$ws.Range("K132").Value = 15949.929
$ws.Range("L132").Value = 14587.2
$ws.Range("M132").Value = -13419.929
$ws.Range("N132").Value = -19647.2

$ws.Range("H136").Value = 1337
$ws.Range("I136").Value = 1204.4
$ws.Range("K136").Value = 3613.2
$ws.Range("M136").Value = -1063.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1412.8
$ws.Range("I113").Value = 1855.5834
$ws.Range("J113").Value = 748.625
$ws.Range("K113").Value = 5566.7502
$ws.Range("L113").Value = 2245.875
$ws.Range("M113").Value = -3396.7502
$ws.Range("N113").Value = -6585.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 11778001
$ws.Range("I11").Value = 13250001
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 13250001
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -13249862
$ws.Range("N11").Value = -2278

$ws.Range("H18").Value = 48602.6
$ws.Range("I18").Value = 16497.5
$ws.Range("K18").Value = 16497.5
$ws.Range("M18").Value = -16204.5

$ws.Range("H53").Value = 21666.666
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 21666.666
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 21666.666
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -22928.666

$ws.Range("H86").Value = 18142
$ws.Range("J86").Value = 18142
$ws.Range("L86").Value = 18142
$ws.Range("N86").Value = -20514

$ws.Range("H89").Value = 18142
$ws.Range("J89").Value = 18142
$ws.Range("L89").Value = 54426
$ws.Range("N89").Value = -66282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 26762
$ws.Range("J57").Value = 26762
$ws.Range("L57").Value = 26762
$ws.Range("N57").Value = -27894

$ws.Range("H82").Value = 2966.7778
$ws.Range("I82").Value = 2702
$ws.Range("J82").Value = 2982.353
$ws.Range("K82").Value = 2702
$ws.Range("L82").Value = 2982.353
$ws.Range("M82").Value = -2341
$ws.Range("N82").Value = -3704.353

$ws.Range("H85").Value = 2966.7778
$ws.Range("I85").Value = 2702
$ws.Range("J85").Value = 2982.353
$ws.Range("K85").Value = 2702
$ws.Range("L85").Value = 2982.353
$ws.Range("M85").Value = -1454
$ws.Range("N85").Value = -5478.353

$ws.Range("H93").Value = 1666
$ws.Range("I93").Value = 1674.1666
$ws.Range("J93").Value = 1633.3334
$ws.Range("K93").Value = 1674.1666
$ws.Range("L93").Value = 1633.3334
$ws.Range("M93").Value = -426.1666
$ws.Range("N93").Value = -4129.3334

$ws.Range("H133").Value = 32958.445
$ws.Range("J133").Value = 32958.445
$ws.Range("L133").Value = 32958.445
$ws.Range("N133").Value = -38018.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 38421.332
$ws.Range("J46").Value = 38421.332
$ws.Range("L46").Value = 38421.332
$ws.Range("N46").Value = -38883.332

$ws.Range("H134").Value = 38421.332
$ws.Range("J134").Value = 38421.332
$ws.Range("L134").Value = 115263.996
$ws.Range("N134").Value = -120333.996
